$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Numeric data refresh (Casos totales / Nuevos casos / Casos activos / Recuperados / Casos criticos / Muertes hoy / Muertes) ---

# Estados Unidos (row 6)
$ws.Range("B6").Value = 46145
$ws.Range("C6").Value = 2411
$ws.Range("E6").Value = 45268
$ws.Range("G6").Value = 29
$ws.Range("H6").Value = 582

# India (row 46)
$ws.Range("D46").Value = 37
$ws.Range("E46").Value = 452

# Row 58 becomes Colombia's updated numbers, row 59 keeps Argentina's (unchanged) numbers
# since Colombia's total (306) overtook Argentina's (301).
$ws.Range("A58").Value = "Colombia"
$ws.Range("B58").Value = 306
$ws.Range("C58").Value = 29
$ws.Range("D58").Value = 6
$ws.Range("E58").Value = 297
$ws.Range("H58").Value = 3

$ws.Range("A59").Value = "Argentina"
$ws.Range("B59").Value = 301
$ws.Range("D59").Value = 51
$ws.Range("E59").Value = 246
$ws.Range("H59").Value = 4

# Kirguistan / Kenia tie-break swap (row 136/137) - identical data, just reorder labels
$ws.Range("A136").Value = "Kenia"
$ws.Range("A137").Value = "Kirguistan"

# Benin / Bermudas tie-break swap (row 151/153) - identical data, just reorder labels
$ws.Range("A151").Value = "Bermudas"
$ws.Range("A153").Value = "Benin"

# Namibia's new case (row 161) triggers resorting of the tied B=4/B=3 block (rows 158-169)
$ws.Range("C161").Value = 1

# Curazao's stats shift onto the new sorted position (row 164)
$ws.Range("E164").Value = 4
$ws.Range("H164").Value = 0

# San Bartolome's stats shift onto the new sorted position (row 165)
$ws.Range("B165").Value = 4
$ws.Range("H165").Value = 1

# Relabel rows 158-169 to match the new sorted order of tied countries
$ws.Range("A158").Value = "Suazilandia"
$ws.Range("A160").Value = "Groenlandia"
$ws.Range("A161").Value = "Fiyi"
$ws.Range("A162").Value = "Guinea"
$ws.Range("A163").Value = "Namibia"
$ws.Range("A164").Value = "Bahamas"
$ws.Range("A165").Value = "Curazao"
$ws.Range("A166").Value = "San Bartolome"
$ws.Range("A167").Value = "Republica de Yibuti"
$ws.Range("A168").Value = "Republica de Africa Central"
$ws.Range("A169").Value = "Zambia"

# --- Timestamp update ---
$ws.Range("A1").Value = "Datos actualizados a 24 de Marzo de 2020 a las 04:46"
